$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reclassified the problem as 3-class and reran setfit: rows are now sorted
# alphabetically by model name, a new "setfit" row is inserted, and the SVM
# numbers (now in row 8) were refreshed.
$data = @(
    @("ComplementNB", 61.42, 59.44, 59.03, 61.42, 60.21871781349182, 5.018226484457652),
    @("Decision Tree", 54.16, 52.66, 51.61, 54.16, 1364.714178800583, 113.7261815667152),
    @("LR", 64.69, 63.02, 62.14, 64.69, 7841.706496715546, 653.4755413929621),
    @("MultinomialNB", 62.81, 57.04, 57.5, 62.81, 59.53350234031677, 4.961125195026398),
    @("RF", 59.41, 46.82, 53.99, 59.41, 5079.660954475403, 423.3050795396169),
    @("setfit", 63.89, 65.29000000000001, 70.23999999999999, 63.89, 1915.019298315048, 159.584941526254),
    @("SVM", 64.62, 62.42, 61.68, 64.62, 438370.2274112701, 54796.27842640877)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $rowIndex++
}

# Row 8 is new (the sheet previously only spanned through row 7); copy the
# label-cell style (bold/border/centered) used throughout column A so the
# new row's "A8" cell matches its siblings.
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null

$ws.Range("A1").Select() | Out-Null
